$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "esboço de bind function na página de vinculação"
# Insert a new field row above row 25 ("DEMANDA-ESPECÍFICA"), pushing it and
# everything below it down by one row, to add the new
# "DEMANDA-COMPLEXIDADE" field between "DEMANDA-ESPECIALIDADE" (B24) and
# "DEMANDA-ESPECÍFICA" (now B26).
$ws.Rows.Item(25).Insert()

$ws.Range("B25").Value = "DEMANDA-COMPLEXIDADE"

# Match the formatting (fill/border) of the neighbouring field cells.
$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)

$ws.Range("B26").Select()
